# Weekly update: insert a new price record (row 243) for
# "Vega Modelo de Temuco - Granada", pushing the existing rows
# 243-276 down to 244-277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 243 (shifts 243:276 -> 244:277,
# carrying their formatting/data down automatically, same as Excel's
# Rows(...).Insert()).
$ws.Rows("243").Insert()

# Populate the newly inserted row 243 with this week's data.
$ws.Range("A243").Value = 10
$ws.Range("B243").Value = "Vega Modelo de Temuco"
$ws.Range("C243").Value = "La Araucanía"
$ws.Range("D243").Value = 45142
$ws.Range("E243").Value = 9
$ws.Range("F243").Value = "Fruta"
$ws.Range("G243").Value = 100104
$ws.Range("H243").Value = "Frutos de pepita"
$ws.Range("I243").Value = 100104001
$ws.Range("J243").Value = "Granada"
$ws.Range("K243").Value = "Wonderfull"
$ws.Range("L243").Value = "Primera"
$ws.Range("M243").Value = 50
$ws.Range("N243").Value = 16000
$ws.Range("O243").Value = 16000
$ws.Range("P243").Value = 16000
$ws.Range("Q243").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R243").Value = "Región de O'Higgins"
$ws.Range("S243").Value = 1600
$ws.Range("T243").Value = 10
